$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 5, column D (Correspond Handoff Datetime) and G (Correspond Handback DateTime)
$wsZh.Range("D5").Value = "2016-01-27 08:29:56"
$wsZh.Range("G5").Value = "2016-01-27 08:30:54"

# de-de sheet: row 5, column D and G
$wsDe.Range("D5").Value = "2016-01-27 08:30:11"
$wsDe.Range("G5").Value = "2016-01-27 08:31:17"
